$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1806451612903226
$ws.Range("C2").Value = 0.5612903225806452
$ws.Range("J2").Value = 0.00967741935483871
$ws.Range("P2").Value = 0.132258064516129
$ws.Range("S2").Value = 0.1161290322580645
$ws.Range("C3").Value = 0.02824858757062147
$ws.Range("J3").Value = 0.01694915254237288
$ws.Range("P3").Value = 0.7344632768361582
$ws.Range("S3").Value = 0.2203389830508475
$ws.Range("J4").Value = 0.0392156862745098
$ws.Range("P4").Value = 0.6274509803921569
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.06854838709677419
$ws.Range("D6").Value = 0.008064516129032258
$ws.Range("F6").Value = 0.1008064516129032
$ws.Range("J6").Value = 0.2056451612903226
$ws.Range("O6").Value = 0.03225806451612903
$ws.Range("Q6").Value = 0.1814516129032258
$ws.Range("R6").Value = 0.07661290322580645
$ws.Range("S6").Value = 0.3266129032258064
$ws.Range("B7").Value = 0.1153846153846154
$ws.Range("D7").Value = 0.02884615384615385
$ws.Range("F7").Value = 0.07692307692307693
$ws.Range("J7").Value = 0.1057692307692308
$ws.Range("O7").Value = 0.02884615384615385
$ws.Range("Q7").Value = 0.2259615384615385
$ws.Range("R7").Value = 0.05288461538461538
$ws.Range("S7").Value = 0.3653846153846154
$ws.Range("B8").Value = 0.1095890410958904
$ws.Range("D8").Value = 0.0182648401826484
$ws.Range("F8").Value = 0.0547945205479452
$ws.Range("J8").Value = 0.091324200913242
$ws.Range("O8").Value = 0.0273972602739726
$ws.Range("Q8").Value = 0.2214611872146119
$ws.Range("R8").Value = 0.06164383561643835
$ws.Range("S8").Value = 0.4155251141552511
$ws.Range("B9").Value = 0.1238532110091743
$ws.Range("D9").Value = 0.03669724770642202
$ws.Range("F9").Value = 0.1192660550458716
$ws.Range("J9").Value = 0.0871559633027523
$ws.Range("O9").Value = 0.004587155963302753
$ws.Range("Q9").Value = 0.1972477064220184
$ws.Range("R9").Value = 0.07798165137614679
$ws.Range("S9").Value = 0.3532110091743119
$ws.Range("B10").Value = 0.1077302631578947
$ws.Range("D10").Value = 0.0287828947368421
$ws.Range("E10").Value = 0.0008223684210526315
$ws.Range("F10").Value = 0.078125
$ws.Range("J10").Value = 0.1077302631578947
$ws.Range("O10").Value = 0.02384868421052632
$ws.Range("Q10").Value = 0.21875
$ws.Range("R10").Value = 0.05921052631578947
$ws.Range("S10").Value = 0.375
$ws.Range("G11").Value = 0.09731543624161074
$ws.Range("J11").Value = 0.07046979865771812
$ws.Range("K11").Value = 0.1644295302013423
$ws.Range("L11").Value = 0.6543624161073825
$ws.Range("S11").Value = 0.01342281879194631
$ws.Range("G12").Value = 0.7644230769230769
$ws.Range("J12").Value = 0.1538461538461539
$ws.Range("L12").Value = 0.04807692307692308
$ws.Range("S12").Value = 0.03365384615384615
$ws.Range("G13").Value = 0.5945945945945946
$ws.Range("J13").Value = 0.3513513513513514
$ws.Range("S13").Value = 0.05405405405405406
$ws.Range("F15").Value = 0.01606425702811245
$ws.Range("H15").Value = 0.1847389558232932
$ws.Range("I15").Value = 0.07630522088353414
$ws.Range("J15").Value = 0.3132530120481928
$ws.Range("K15").Value = 0.05220883534136546
$ws.Range("M15").Value = 0.004016064257028112
$ws.Range("O15").Value = 0.06024096385542169
$ws.Range("S15").Value = 0.2931726907630522
$ws.Range("F16").Value = 0.02040816326530612
$ws.Range("H16").Value = 0.1683673469387755
$ws.Range("I16").Value = 0.0663265306122449
$ws.Range("J16").Value = 0.3979591836734694
$ws.Range("K16").Value = 0.1377551020408163
$ws.Range("M16").Value = 0.00510204081632653
$ws.Range("N16").Value = 0.00510204081632653
$ws.Range("O16").Value = 0.08673469387755102
$ws.Range("S16").Value = 0.1122448979591837
$ws.Range("F17").Value = 0.01217038539553753
$ws.Range("H17").Value = 0.1724137931034483
$ws.Range("I17").Value = 0.1135902636916836
$ws.Range("J17").Value = 0.4279918864097363
$ws.Range("K17").Value = 0.1075050709939148
$ws.Range("M17").Value = 0.01825557809330629
$ws.Range("N17").Value = 0.008113590263691683
$ws.Range("O17").Value = 0.05882352941176471
$ws.Range("S17").Value = 0.08113590263691683
$ws.Range("F18").Value = 0.006896551724137931
$ws.Range("H18").Value = 0.1793103448275862
$ws.Range("I18").Value = 0.07586206896551724
$ws.Range("J18").Value = 0.4275862068965517
$ws.Range("K18").Value = 0.1103448275862069
$ws.Range("M18").Value = 0.006896551724137931
$ws.Range("O18").Value = 0.09655172413793103
$ws.Range("S18").Value = 0.09655172413793103
$ws.Range("F19").Value = 0.01869918699186992
$ws.Range("H19").Value = 0.2056910569105691
$ws.Range("I19").Value = 0.09674796747967479
$ws.Range("J19").Value = 0.3772357723577236
$ws.Range("K19").Value = 0.1121951219512195
$ws.Range("M19").Value = 0.02113821138211382
$ws.Range("N19").Value = 0.0008130081300813008
$ws.Range("O19").Value = 0.06829268292682927
$ws.Range("S19").Value = 0.09918699186991869
